$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map, derived from the day-over-day cryptos data refresh
$updates = @(
    @{ Cell = 'D2'; Value = '30.046.20' }
    @{ Cell = 'E2'; Value = '  -0.73%  ' }
    @{ Cell = 'D3'; Value = '1.916.62' }
    @{ Cell = 'E3'; Value = '  +0.28%  ' }
    @{ Cell = 'D4'; Value = '1.001' }
    @{ Cell = 'E4'; Value = '  +0.01%  ' }
    @{ Cell = 'D5'; Value = '320.23' }
    @{ Cell = 'E5'; Value = '  -2.68%  ' }
    @{ Cell = 'E6'; Value = '  +0.04%  ' }
    @{ Cell = 'D7'; Value = '0.5026' }
    @{ Cell = 'E7'; Value = '  -3.08%  ' }
    @{ Cell = 'D8'; Value = '0.4033' }
    @{ Cell = 'E8'; Value = '  -0.18%  ' }
    @{ Cell = 'D9'; Value = '0.08236' }
    @{ Cell = 'E9'; Value = '  -3.08%  ' }
    @{ Cell = 'D10'; Value = '1.108' }
    @{ Cell = 'E10'; Value = '  -1.50%  ' }
    @{ Cell = 'D11'; Value = '42.03' }
    @{ Cell = 'E11'; Value = '  -1.83%  ' }
    @{ Cell = 'D12'; Value = '23.75' }
    @{ Cell = 'E12'; Value = '  +0.39%  ' }
    @{ Cell = 'B13'; Value = 'Polkadot' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' }
    @{ Cell = 'D13'; Value = '6.415' }
    @{ Cell = 'E13'; Value = '  -0.31%  ' }
    @{ Cell = 'B14'; Value = 'WrappedEther' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Cell = 'D14'; Value = '1.910.39' }
    @{ Cell = 'E14'; Value = '  -0.57%  ' }
    @{ Cell = 'D15'; Value = '7.289' }
    @{ Cell = 'E15'; Value = '  -1.32%  ' }
    @{ Cell = 'D16'; Value = '1.002' }
    @{ Cell = 'E16'; Value = '  +0.09%  ' }
    @{ Cell = 'D17'; Value = '92.24' }
    @{ Cell = 'E17'; Value = '  -3.07%  ' }
    @{ Cell = 'D18'; Value = '0.00001097' }
    @{ Cell = 'E18'; Value = '  -1.62%  ' }
    @{ Cell = 'D19'; Value = '0.06502' }
    @{ Cell = 'E19'; Value = '  -3.21%  ' }
    @{ Cell = 'D20'; Value = '18.18' }
    @{ Cell = 'E20'; Value = '  -1.56%  ' }
    @{ Cell = 'E21'; Value = '  -0.05%  ' }
    @{ Cell = 'D22'; Value = '5.938' }
    @{ Cell = 'E22'; Value = '  -1.31%  ' }
    @{ Cell = 'D23'; Value = '30.085.35' }
    @{ Cell = 'E23'; Value = '  -0.63%  ' }
    @{ Cell = 'E24'; Value = '  -0.88%  ' }
    @{ Cell = 'E25'; Value = '  -1.27%  ' }
    @{ Cell = 'D26'; Value = '22.33' }
    @{ Cell = 'E26'; Value = '  +1.76%  ' }
    @{ Cell = 'D27'; Value = '2.133.27' }
    @{ Cell = 'E27'; Value = '  -0.33%  ' }
    @{ Cell = 'D28'; Value = '162.13' }
    @{ Cell = 'E28'; Value = '  -0.22%  ' }
    @{ Cell = 'D29'; Value = '2.290' }
    @{ Cell = 'E29'; Value = '  -4.76%  ' }
    @{ Cell = 'E30'; Value = '  -0.57%  ' }
    @{ Cell = 'D31'; Value = '1.133' }
    @{ Cell = 'E31'; Value = '  +2.32%  ' }
    @{ Cell = 'D32'; Value = '0.1038' }
    @{ Cell = 'E32'; Value = '  -2.41%  ' }
    @{ Cell = 'D33'; Value = '6.020' }
    @{ Cell = 'E33'; Value = '  +0.08%  ' }
    @{ Cell = 'D34'; Value = '3.812' }
    @{ Cell = 'E34'; Value = '  +4.57%  ' }
    @{ Cell = 'D35'; Value = '0.02442' }
    @{ Cell = 'E35'; Value = '  -2.04%  ' }
    @{ Cell = 'D36'; Value = '5.351' }
    @{ Cell = 'E36'; Value = '  +2.97%  ' }
    @{ Cell = 'D37'; Value = '0.06435' }
    @{ Cell = 'E37'; Value = '  -2.31%  ' }
    @{ Cell = 'B38'; Value = 'Algorand' }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' }
    @{ Cell = 'D38'; Value = '0.2165' }
    @{ Cell = 'E38'; Value = '  -2.17%  ' }
    @{ Cell = 'B39'; Value = 'FraxShare' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Cell = 'D39'; Value = '8.907' }
    @{ Cell = 'E39'; Value = '  +1.14%  ' }
    @{ Cell = 'D40'; Value = '1.206' }
    @{ Cell = 'E40'; Value = '  -2.26%  ' }
    @{ Cell = 'D41'; Value = '0.6426' }
    @{ Cell = 'E41'; Value = '  -1.64%  ' }
    @{ Cell = 'D42'; Value = '11.40' }
    @{ Cell = 'E42'; Value = '  -4.83%  ' }
    @{ Cell = 'D43'; Value = '1.221' }
    @{ Cell = 'E43'; Value = '  -1.29%  ' }
    @{ Cell = 'E44'; Value = '  +0.01%  ' }
    @{ Cell = 'D45'; Value = '13.40' }
    @{ Cell = 'E45'; Value = '  +0.71%  ' }
    @{ Cell = 'D46'; Value = '2.193' }
    @{ Cell = 'E46'; Value = '  +5.84%  ' }
    @{ Cell = 'D47'; Value = '0.5997' }
    @{ Cell = 'E47'; Value = '  -2.40%  ' }
    @{ Cell = 'E48'; Value = '  -2.67%  ' }
    @{ Cell = 'D49'; Value = '123.49' }
    @{ Cell = 'E49'; Value = '  -1.56%  ' }
    @{ Cell = 'D50'; Value = '1.215' }
    @{ Cell = 'E50'; Value = '  -2.49%  ' }
    @{ Cell = 'D51'; Value = '78.89' }
    @{ Cell = 'E51'; Value = '  -0.79%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "1.001") are not
    # coerced into numbers, matching the original inline-string cell type.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
